$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 4 - Baseline Adj R^2 (label unchanged)
$ws.Range("B4").Value = -0.14871979361862
$ws.Range("C4").Value = -0.1551244756811565
$ws.Range("D4").Value = -0.2282119191780376
$ws.Range("E4").Value = -0.1129285212772451
$ws.Range("F4").Value = -0.2658332717264276
$ws.Range("G4").Value = 0.005122818527443562

# Row 7 - Model R^2 (label unchanged)
$ws.Range("B7").Value = 0.4478752838460561
$ws.Range("C7").Value = 0.4346435041987016
$ws.Range("D7").Value = 0.5641682672869548
$ws.Range("E7").Value = 0.6330952554014504
$ws.Range("F7").Value = 0.4476816913065088
$ws.Range("G7").Value = 0.4974829832875998

# Row 8 - Model Adj R^2 (label unchanged)
$ws.Range("B8").Value = 0.4459509056731166
$ws.Range("C8").Value = 0.431739985878372
$ws.Range("D8").Value = 0.5022226403023595
$ws.Range("E8").Value = 0.4774386970869142
$ws.Range("F8").Value = 0.4387294987488141
$ws.Range("G8").Value = 0.4679851290571024

# Row 9 - Model RMSE (label unchanged)
$ws.Range("B9").Value = 2.340395714487467
$ws.Range("C9").Value = 2.37547511409402
$ws.Range("D9").Value = 1.759236457395349
$ws.Range("E9").Value = 1.850913903771682
$ws.Range("F9").Value = 2.256300916109789
$ws.Range("G9").Value = 2.544226272518651

# Row 10 - Model HH (label unchanged)
$ws.Range("B10").Value = 275
$ws.Range("C10").Value = 190
$ws.Range("D10").Value = 6
$ws.Range("E10").Value = 5
$ws.Range("F10").Value = 48
$ws.Range("G10").Value = 26

# Row 11 - label "Lift R²" -> "Delta R²"
$ws.Range("A11").Value = "Delta R²"
$ws.Range("B11").Value = 0.5926052301383323
$ws.Range("C11").Value = 0.5838658840187134
$ws.Range("D11").Value = 0.6395360365228366
$ws.Range("E11").Value = 0.4145131533195161
$ws.Range("F11").Value = 0.6933250845889946
$ws.Range("G11").Value = 0.4371986576686085

# Row 12 - label "Lift Adj R²" -> "Delta Adj R²"
$ws.Range("A12").Value = "Delta Adj R²"
$ws.Range("B12").Value = 0.5946706992917365
$ws.Range("C12").Value = 0.5868644615595285
$ws.Range("D12").Value = 0.7304345594803971
$ws.Range("E12").Value = 0.5903672183641593
$ws.Range("F12").Value = 0.7045627704752416
$ws.Range("G12").Value = 0.4628623105296589

# Row 13 - label "Drop RMSE" -> "Delta RMSE"
$ws.Range("A13").Value = "Delta RMSE"
$ws.Range("B13").Value = -1.029544698560945
$ws.Range("C13").Value = -1.011338752257862
$ws.Range("D13").Value = -1.004160661435681
$ws.Range("E13").Value = -0.8502518253712559
$ws.Range("F13").Value = -1.132134124928308
$ws.Range("G13").Value = -0.9349672172528005

# Row 14 - label "Drop HH" -> "Delta HH"
$ws.Range("A14").Value = "Delta HH"
$ws.Range("B14").Value = 35
$ws.Range("C14").Value = 25
$ws.Range("D14").Value = 2
$ws.Range("E14").Value = 0
$ws.Range("F14").Value = 3
$ws.Range("G14").Value = 5
